# Add files via upload
# - Fix a duplicated path segment ("16x16_modified16x16" -> "16x16_modified")
#   in the "background" image-path rows (A29:A58).
# - Update the sheet view: scroll so row 19 is at the top, and select
#   A29:A58 (active cell A29) to match where the data was edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 29; $r -le 58; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($old -ne $null -and $old -like "*16x16_modified16x16*") {
        $new = $old -replace "16x16_modified16x16", "16x16_modified"
        $cell.Value = $new
    }
}

# Reproduce the new selection (A29:A58, active cell A29) and scroll position
# (top-left visible row 19) recorded in the saved sheetView.
$ws.Range("A29:A58").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
